# Add two new columns, I ("I0") and J ("IF"), to the sheet - mirrors the
# existing H ("IP") column's header formatting and fills in the data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers (row 1): copy the style from the existing H1 header cell so
# the new header cells pick up the same bold/border/center formatting,
# then overwrite the text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows 2-34 for columns I and J.
$iValues = @(2,7,8,4,9,9,4,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,4,1,1,2,1,1)
$jValues = @(3,7,9,4,9,9,6,3,5,2,6,6,6,5,7,7,6,6,5,5,7,5,6,6,7,6,6,8,5,4,4,3,2)

for ($i = 0; $i -lt $iValues.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 9).Value = $iValues[$i]
    $ws.Cells.Item($row, 10).Value = $jValues[$i]
}
